# Updated code to use ECI frame for thrust vector
# Add new ECI X/Y/Z Thrust Component columns (K, L, M) alongside the
# existing RIC Thrust Component columns (H, I, J), with header labels and
# data for the two burn rows (row 4 and row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (row 1)
$ws.Range("K1").Value = "ECI X Thrust Component"
$ws.Range("L1").Value = "ECI Y Thrust Component"
$ws.Range("M1").Value = "ECI Z Thrust Component"

# Burn1 segment (row 4) ECI thrust components
$ws.Range("K4").Value = -0.5083
$ws.Range("L4").Value = 0.09305
$ws.Range("M4").Value = -0.8561

# Burn2 segment (row 6) ECI thrust components
$ws.Range("K6").Value = -0.38537
$ws.Range("L6").Value = -0.18610758
$ws.Range("M6").Value = 0.903796

# Update the saved selection / view to match the edited workbook
$ws.Range("I11").Select()
